$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# "Contenu du stage" breakdown (rows 16-20 change; 21-23 stay the same)
$ws.Range("E16").Value = 1
Set-TextValue "G16" "2.78 %"

$ws.Range("E17").Value = 33
Set-TextValue "G17" "91.67 %"

$ws.Range("E19").Value = 1
Set-TextValue "G19" "2.78 %"

$ws.Range("E20").Value = 1
Set-TextValue "G20" "2.78 %"

# Nudge the "Contenu du stage" pie chart to pick up the refreshed
# D16:E23 range (best effort - harmless if the host ignores it).
try {
    $chartObj = $ws.ChartObjects().Item(2)
    $chartObj.Chart.SetSourceData($ws.Range("D16:E23"))
    $chartObj.Chart.Refresh()
} catch {
}
